$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "List of Projects": insert a new "Type" column, add the finished
# "Reverse a string" project row, and tidy up the two punctuation-only
# description edits.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("List of Projects")

# Insert a new column D ("Type") between Description and Date Started.
$ws1.Columns.Item(4).Insert()
$ws1.Columns.Item(4).ColumnWidth = 10.2

# Project header label is now left aligned instead of centered.
$ws1.Range("B3").HorizontalAlignment = -4131
$ws1.Range("B3").VerticalAlignment = -4108

# New finished project: "Reverse a string" (entered first, as it was in the
# live editing session).
$ws1.Range("B6").Value = "Reverse a string"
$ws1.Range("B6").VerticalAlignment = -4108

# Header for the new column.
$ws1.Range("D3").Value = "Type"
$ws1.Range("D3").HorizontalAlignment = -4108
$ws1.Range("D3").VerticalAlignment = -4108
$ws1.Range("D3").WrapText = $true
$ws1.Range("D3").Font.Bold = $true

# Type values for the existing two rows plus the new one.
$ws1.Range("D4").Value = "Numbers"
$ws1.Range("D4").HorizontalAlignment = -4108
$ws1.Range("D4").VerticalAlignment = -4108
$ws1.Range("D4").WrapText = $true

$ws1.Range("D5").Value = "Numbers"
$ws1.Range("D5").HorizontalAlignment = -4108
$ws1.Range("D5").VerticalAlignment = -4108
$ws1.Range("D5").WrapText = $true

$ws1.Range("D6").Value = "Text"
$ws1.Range("D6").HorizontalAlignment = -4108
$ws1.Range("D6").VerticalAlignment = -4108
$ws1.Range("D6").WrapText = $true

# Description of the new row.
$ws1.Range("C6").Value = "enter a string and the program reverses the string and prints it out."
$ws1.Range("C6").VerticalAlignment = -4108
$ws1.Range("C6").WrapText = $true

# Finished descriptions gained a trailing period.
$ws1.Range("C5").Value = "enter a number and have the program generate the fibonacci sequence to that number or Nth number."
$ws1.Range("C4").Value = "write a program to determine the Nth digit of PI."

# Dates now centre horizontally as well as vertically.
$ws1.Range("E4:F5").HorizontalAlignment = -4108
$ws1.Range("E4:F5").VerticalAlignment = -4108

$ws1.Range("E6").Value = 43516
$ws1.Range("F6").Value = 43516
$ws1.Range("E6:F6").NumberFormat = "m/d/yy"
$ws1.Range("E6:F6").HorizontalAlignment = -4108
$ws1.Range("E6:F6").VerticalAlignment = -4108

$ws1.Range("G6").Value = "Java"
$ws1.Range("G6").HorizontalAlignment = -4108
$ws1.Range("G6").VerticalAlignment = -4108

$ws1.Rows.Item(6).RowHeight = 29

# ---------------------------------------------------------------------------
# Sheet "Sources": link the Martyr2 mega-list URL cell to its address.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sources")
$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://github.com/JSmolins/Martyrs-Mega-List") | Out-Null

# ---------------------------------------------------------------------------
# Restore the on-screen selections recorded in the saved workbook.
# ---------------------------------------------------------------------------
$ws2.Range("C9").Select() | Out-Null
$ws1.Range("D10").Select() | Out-Null
